# Add the final prep-work row (row 15) for the new patient "Mateusz Lugowski"
# to the visits sheet. Columns: First name | Last name | phone | visit date |
# visit time | e-mail | queue number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C, E and G hold values that look numeric ("11111111111", "10.00",
# "15") but must be stored as text, like every other such column in this
# sheet (phone numbers, times, queue numbers are all text). Mark them as
# text before writing so Excel keeps them as strings instead of coercing
# them to numbers.
$ws.Range("C15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("G15").NumberFormat = "@"

$ws.Range("A15").Value = "Mateusz"
$ws.Range("B15").Value = "Lugowski"
$ws.Range("C15").Value = "11111111111"
$ws.Range("D15").Value = "19.01.2023"
$ws.Range("E15").Value = "10.00"
$ws.Range("F15").Value = "lugowski.mateusz.02@gmail.com"
$ws.Range("G15").Value = "15"
